# Přidána historie chyb kvadratického průměru do algoritmu zpětné propagace
# Update a few probability values in the diagnosis output table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0
$ws.Range("F3").Value = 0.01
$ws.Range("H3").Value = 0.01
$ws.Range("H5").Value = 0.01
$ws.Range("G6").Value = 0
